$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a new empty paragraph (underlined "Pericles" paragraph mark) at the
#    very start of the document, before the existing "General Skills" heading.
# ---------------------------------------------------------------------------
$introRange = $d.Range(0, 0)
$introRange.InsertXML('<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Pericles" w:hAnsi="Pericles"/><w:u w:val="single"/></w:rPr></w:pPr></w:p>')

# ---------------------------------------------------------------------------
# 2. "Thief Skills" heading paragraph loses its <w:lastRenderedPageBreak/>.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Thief Skills", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Paragraphs(1).Range.InsertXML('<w:p w14:paraId="55CEB32F" w14:textId="77777777" w:rsidR="000D1198" w:rsidRDefault="000D1198"><w:pPr><w:rPr><w:rFonts w:ascii="Pericles" w:hAnsi="Pericles"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Pericles" w:hAnsi="Pericles"/></w:rPr><w:t>Thief Skills</w:t></w:r></w:p>')

# ---------------------------------------------------------------------------
# 3. "Scribe" table-cell paragraph gains a <w:lastRenderedPageBreak/>.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Scribe", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Paragraphs(1).Range.InsertXML('<w:p w14:paraId="44949363" w14:textId="77777777" w:rsidR="0080396E" w:rsidRDefault="0080396E" w:rsidP="0080396E"><w:pPr><w:rPr><w:rFonts w:ascii="Pericles" w:hAnsi="Pericles"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Pericles" w:hAnsi="Pericles"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:lastRenderedPageBreak/><w:t>Scribe</w:t></w:r></w:p>')

# ---------------------------------------------------------------------------
# 4. "License: Forgery (...)" bullet loses its <w:lastRenderedPageBreak/>.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("License: Forgery (can create temporary fake credentials for various purposes)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Paragraphs(1).Range.InsertXML('<w:p w14:paraId="2C0A2761" w14:textId="0B0CD391" w:rsidR="00591302" w:rsidRDefault="00591302" w:rsidP="0080396E"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Pericles" w:hAnsi="Pericles"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Pericles" w:hAnsi="Pericles"/><w:b/><w:bCs/><w:i/><w:iCs/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>License: Forgery (can create temporary fake credentials for various purposes)</w:t></w:r><w:r w:rsidR="000C7C83"><w:rPr><w:rFonts w:ascii="Pericles" w:hAnsi="Pericles"/><w:b/><w:bCs/><w:i/><w:iCs/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="000C7C83"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:i/><w:iCs/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>–</w:t></w:r><w:r w:rsidR="000C7C83"><w:rPr><w:rFonts w:ascii="Pericles" w:hAnsi="Pericles"/><w:b/><w:bCs/><w:i/><w:iCs/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> 10000 G</w:t></w:r></w:p>')

# ---------------------------------------------------------------------------
# 5. "Tailor" table-cell paragraph loses its <w:lastRenderedPageBreak/>.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Tailor", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Paragraphs(1).Range.InsertXML('<w:p w14:paraId="2E61BEFE" w14:textId="77777777" w:rsidR="0080396E" w:rsidRDefault="0080396E" w:rsidP="0080396E"><w:pPr><w:rPr><w:rFonts w:ascii="Pericles" w:hAnsi="Pericles"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Pericles" w:hAnsi="Pericles"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>Tailor</w:t></w:r></w:p>')

# ---------------------------------------------------------------------------
# 6. Split the "features defensive moves..." run, adding gramStart/gramEnd
#    <w:proofErr/> markers around the word "features" (spear & shield row).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("features defensive moves, impaling attacks and taunts, and spear throwing", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Paragraphs(1).Range.InsertXML('<w:p w14:paraId="580914C9" w14:textId="77777777" w:rsidR="00196A61" w:rsidRDefault="00196A61" w:rsidP="00723274"><w:r><w:rPr><w:rFonts w:ascii="Pericles" w:hAnsi="Pericles"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">Spear and shield </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>–</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Pericles" w:hAnsi="Pericles"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Pericles" w:hAnsi="Pericles"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>features</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Pericles" w:hAnsi="Pericles"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> defensive moves, impaling attacks and taunts, and spear throwing</w:t></w:r></w:p>')

Write-Output "done"
